# Update the generic "eBird Observation Dataset" labels on Sheet2 so each
# region (Hawaii / Utah / New Zealand) gets its own descriptive dataset name.
#
# Shared-string insertion order matters for byte-for-byte parity with the
# target workbook: new unique strings are appended to the shared string
# table in the order cells are first written, so we set the rows in the
# same order the author edited them (NZ, then UT, then HI) to reproduce the
# resulting index layout (NZ=48, UT=49, HI=50).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B19").Value = "NZ eBird Observation Dataset"
$ws.Range("B17").Value = "UT eBird Observation Dataset"
$ws.Range("B15").Value = "HI eBird Observation Dataset"

# Reproduce the author's final on-screen selection / scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
